$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '58.824.18'
$ws.Range('E2').Value = '  -0.61%  '
$ws.Range('D3').Value = '2.511.04'
$ws.Range('E3').Value = '  +0.01%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '531.66'
$ws.Range('E5').Value = '  -2.09%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '138.27'
$ws.Range('E6').Value = '  -4.11%  '
$ws.Range('E7').Value = '  +0.43%  '
$ws.Range('E8').Value = '  -2.08%  '
$ws.Range('D9').Value = '2.512.00'
$ws.Range('E9').Value = '  -0.99%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.1000'
$ws.Range('E10').Value = '  -0.95%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.160'
$ws.Range('E11').Value = '  +0.71%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '5.46'
$ws.Range('E12').Value = '  -1.62%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '0.356'
$ws.Range('E13').Value = '  +0.22%  '
$ws.Range('D14').Value = '2.955.76'
$ws.Range('E14').Value = '  +0.04%  '
$ws.Range('E15').Value = '  -2.71%  '
$ws.Range('D16').Value = '58.795.44'
$ws.Range('E16').Value = '  -0.52%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '0.0000139'
$ws.Range('E17').Value = '  -0.61%  '
$ws.Range('D18').Value = '2.510.87'
$ws.Range('E18').Value = '  -0.63%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '11.01'
$ws.Range('E19').Value = '  -2.02%  '
$ws.Range('E20').Value = '  -0.51%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '322.40'
$ws.Range('E21').Value = '  -0.90%  '
$ws.Range('E22').Value = '  +0.22%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '5.79'
$ws.Range('E23').Value = '  -0.35%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '62.22'
$ws.Range('E24').Value = '  +0.58%  '
$ws.Range('E25').Value = '  -3.25%  '
$ws.Range('E26').Value = '  +1.51%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '0.999'
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '7.77'
$ws.Range('E28').Value = '  -2.66%  '
$ws.Range('B29').Value = 'PEPE'
$ws.Range('C29').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D29').Value = '0.0₃0768'
$ws.Range('E29').Value = '  -1.58%  '
$ws.Range('B30').Value = 'Aptos'
$ws.Range('C30').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '6.66'
$ws.Range('E30').Value = '  -0.02%  '
$ws.Range('E31').Value = '  -2.00%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '163.17'
$ws.Range('E32').Value = '  +3.38%  '
$ws.Range('E34').Value = '  -8.09%  '
$ws.Range('E35').Value = '  -4.14%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '18.42'
$ws.Range('E36').Value = '  -1.47%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '4.23'
$ws.Range('E37').Value = '  -3.39%  '
$ws.Range('E38').Value = '  -3.81%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '36.74'
$ws.Range('E39').Value = '  -0.49%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '3.64'
$ws.Range('E40').Value = '  -1.79%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.800'
$ws.Range('E41').Value = '  -2.66%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '5.21'
$ws.Range('E42').Value = '  -8.16%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '278.53'
$ws.Range('E43').Value = '  -6.89%  '
$ws.Range('E44').Value = '  +0.55%  '
$ws.Range('E45').Value = '  +0.74%  '
$ws.Range('E46').Value = '  -1.17%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '0.0932'
$ws.Range('E47').Value = '  +0.13%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '121.53'
$ws.Range('E48').Value = '  -1.29%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '18.41'
$ws.Range('E49').Value = '  -1.97%  '
$ws.Range('E50').Value = '  -1.39%  '
$ws.Range('E51').Value = '  -2.52%  '
